$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$title = "Receipee: distinguish which file is bigger for 2 version of file with same structure "

$body = @"
1. Find out the key set of the file structure, if below 2 sqls are both 0, means the set can be used to key this 2 files:
>    SELECT count(*)                       
>    FROM zusrlib/lnp00701                 
>    GROUP BY lhnote, lhrecn, lhpost        
>    having count(*) > 1;
>    SELECT count(*)                       
>    FROM ieom/lnp00701                 
>    GROUP BY lhnote, lhrecn, lhpost        
>    having count(*) > 1                   
2. Use the key set to calculation how much common record between 2 files:
>    select count(*)                            
>    from ifrs201110/lnp00701                   
>    where (lhnote, lhrecn, lhpost) not in (    
>        select lhnote, lhrecn, lhpost              
>        from ifrs201111/lnp00701                   
>    )  ;
"@

# Trim the single trailing newline introduced by the here-string terminator
$body = $body.TrimEnd("`r", "`n")

$ws.Range("A40").Value = "DB2"
$ws.Range("B40").Value = $title
$ws.Range("C40").Value = $body
$ws.Range("C40").WrapText = $true
$ws.Rows.Item(40).RowHeight = 37.5

[void]$ws.Range("B40").Select()
